# Update "想去人数" (interested-count) figures in both the "展览" sheet
# and the "全部类型" sheet to reflect newly generated stats.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 313
$wsExhibit.Range("F5").Value  = 200
$wsExhibit.Range("F6").Value  = 339
$wsExhibit.Range("F7").Value  = 224
$wsExhibit.Range("F8").Value  = 2175
$wsExhibit.Range("F10").Value = 5348
$wsExhibit.Range("F11").Value = 120

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value  = 313
$wsAll.Range("F6").Value  = 200
$wsAll.Range("F7").Value  = 339
$wsAll.Range("F8").Value  = 224
$wsAll.Range("F11").Value = 2175
$wsAll.Range("F13").Value = 5348
$wsAll.Range("F14").Value = 120
